# Applies the "updated summary charts and summary reports including
# comments from Prof. Erhardt" edit to the FAC summary report workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Sheet1: Year 1 changed from 2005 to 2012 (two places)
# ---------------------------------------------------------------------
$ws1.Range("C1").Value2 = 2012
$ws1.Range("E7").Value2 = 2012

# ---------------------------------------------------------------------
# Sheet1: updated "Average Values" data (columns E/F), the % Diff
# formulas (switched from *100/ to a straight ratio, now formatted as a
# percentage), and the recomputed "Riddership Effect" values (column H)
# together with their associated formulas (column I).
# ---------------------------------------------------------------------

# Row 8 - Vehicle Revenue Miles
$ws1.Range("E8").Value2 = 1323080.114
$ws1.Range("F8").Value2 = 1878220.656
$ws1.Range("G8").Formula = "=IFERROR((F8-E8)/E8,0)"
$ws1.Range("H8").Value2 = 654288.211225
$ws1.Range("I8").Formula = "=IFERROR(H8/$E$21,0)"

# Row 9 - Average Fare (2018$)
$ws1.Range("E9").Value2 = 0.597848936
$ws1.Range("F9").Value2 = 0.529305746
$ws1.Range("G9").Formula = "=IFERROR((F9-E9)/E9,0)"
$ws1.Range("H9").Value2 = 27446.91782800001
$ws1.Range("I9").Formula = "=IFERROR(H9/$E$21,0)"

# Row 10 - Population + Employment
$ws1.Range("E10").Value2 = 700037.67
$ws1.Range("F10").Value2 = 910110.75
$ws1.Range("G10").Formula = "=IFERROR((F10-E10)/E10,0)"
$ws1.Range("H10").Value2 = 309126.8709040001
$ws1.Range("I10").Formula = "=IFERROR(H10/$E$21,0)"

# Row 11 - % of Population in Transit Supportive Density
$ws1.Range("E11").Value2 = 2.389578005
$ws1.Range("F11").Value2 = 1.850524206
$ws1.Range("G11").Formula = "=IFERROR((F11-E11)/E11,0)"
$ws1.Range("H11").Value2 = -5778.235695199999
$ws1.Range("I11").Formula = "=IFERROR(H11/$E$21,0)"

# Row 12 - Average Gas Price (2018$)
$ws1.Range("E12").Value2 = 3.9458
$ws1.Range("F12").Value2 = 2.71
$ws1.Range("G12").Formula = "=IFERROR((F12-E12)/E12,0)"
$ws1.Range("H12").Value2 = -154537.2938875
$ws1.Range("I12").Formula = "=IFERROR(H12/$E$21,0)"

# Row 13 - Median Per Capita (2018$)
$ws1.Range("E13").Value2 = 25662.96
$ws1.Range("F13").Value2 = 28855.25
$ws1.Range("G13").Formula = "=IFERROR((F13-E13)/E13,0)"
$ws1.Range("H13").Value2 = -59409.491439
$ws1.Range("I13").Formula = "=IFERROR(H13/$E$21,0)"

# Row 14 - % of Households with 0 Vehicles
$ws1.Range("E14").Value2 = 6.89
$ws1.Range("F14").Value2 = 5.8
$ws1.Range("G14").Formula = "=IFERROR((F14-E14)/E14,0)"
$ws1.Range("H14").Value2 = -23121.973418
$ws1.Range("I14").Formula = "=IFERROR(H14/$E$21,0)"

# Row 15 - % Working at Home
$ws1.Range("E15").Value2 = 4.2
$ws1.Range("F15").Value2 = 4.9
$ws1.Range("G15").Formula = "=IFERROR((F15-E15)/E15,0)"
$ws1.Range("H15").Value2 = -747.2774631
$ws1.Range("I15").Formula = "=IFERROR(H15/$E$21,0)"

# Row 16 - Years Since Ride-hail Start (E/F stay blank)
$ws1.Range("G16").Formula = "=IFERROR((F16-E16)/E16,0)"
$ws1.Range("H16").Value2 = -352540.41948
$ws1.Range("I16").Formula = "=IFERROR(H16/$E$21,0)"

# Row 17 - Bike Share
$ws1.Range("E17").Value2 = 0
$ws1.Range("F17").Value2 = 1
$ws1.Range("G17").Formula = "=IFERROR((F17-E17)/E17,0)"
$ws1.Range("H17").Value2 = 37265.35403
$ws1.Range("I17").Formula = "=IFERROR(H17/$E$21,0)"

# Row 18 - Electric Scooters
$ws1.Range("E18").Value2 = 0
$ws1.Range("F18").Value2 = 0
$ws1.Range("G18").Formula = "=IFERROR((F18-E18)/E18,0)"
$ws1.Range("H18").Value2 = 0
$ws1.Range("I18").Formula = "=IFERROR(H18/$E$21,0)"

# Row 19 - New Reporters
$ws1.Range("G19").Formula = "=IFERROR((F19-E19)/E19,0)"
$ws1.Range("H19").Value2 = 0
$ws1.Range("I19").Formula = "=IFERROR(H19/$E$21,0)"

# Row 20 - Total Modeled Ridership
$ws1.Range("E20").Value2 = 2632910.865
$ws1.Range("F20").Value2 = 3211828.945
$ws1.Range("G20").Formula = "=IFERROR((F20-E20)/E20,0)"

# Row 21 - Total Observed Ridership
$ws1.Range("E21").Value2 = 2968747.833
$ws1.Range("F21").Value2 = 2133709.942
$ws1.Range("G21").Formula = "=IFERROR((F21-E21)/E21,0)"

# ---------------------------------------------------------------------
# Number formats: % Diff columns (G, I) now render as percentages, and
# the Average Values / Riddership Effect columns (E, F, H) use a
# fixed two-decimal numeric format.
# ---------------------------------------------------------------------
$ws1.Range("E8:F21").NumberFormat = "0.00"
$ws1.Range("H8:H21").NumberFormat = "0.00"
$ws1.Range("G8:G21").NumberFormat = "0.00%"
$ws1.Range("I8:I21").NumberFormat = "0.00%"

# ---------------------------------------------------------------------
# Sheet1: sheet view - scroll reset, selection moved
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("H21").Select()

# ---------------------------------------------------------------------
# Sheet2: Average Values header style refresh (cosmetic format nudge)
# ---------------------------------------------------------------------
$ws2.Range("E5").NumberFormat = $ws2.Range("E5").NumberFormat
